# "starter region 1 results"
#
# Fills in the next block of VRP model-rerun results (rows 226-257,
# columns H:L) on the "model_rerun_results_tracker" sheet, and updates
# the saved selection/scroll state on the sheets that were being looked
# at while the data was entered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Enter the new result rows on model_rerun_results_tracker
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("model_rerun_results_tracker")
$ws.Activate()

# Columns: row, H (vehicles), I (distance), J (duration)
# K and L are both the "T" completion flag.
$results = @(
    @(226, 6, 2586, 1335.41),
    @(227, 6, 2612, 1360.19),
    @(228, 6, 2613, 1362.96),
    @(229, 6, 2478, 1279.94),
    @(230, 4, 2128, 1051.06),
    @(231, 5, 2284, 1169.03),
    @(232, 5, 2341, 1202.45),
    @(233, 5, 2329, 1182.86),
    @(234, 6, 2586, 1335.41),
    @(235, 6, 2612, 1360.19),
    @(236, 6, 2613, 1362.96),
    @(237, 6, 2478, 1279.94),
    @(238, 4, 2135, 1063.66),
    @(239, 5, 2284, 1169.03),
    @(240, 6, 2558, 1327.72),
    @(241, 4, 2218, 1143.44),
    @(242, 6, 2586, 1335.41),
    @(243, 6, 2612, 1360.19),
    @(244, 6, 2613, 1362.96),
    @(245, 6, 2478, 1279.94),
    @(246, 4, 2135, 1063.66),
    @(247, 5, 2284, 1169.03),
    @(248, 6, 2558, 1327.72),
    @(249, 4, 2218, 1143.44),
    @(250, 6, 2586, 1335.41),
    @(251, 6, 2612, 1360.19),
    @(252, 6, 2613, 1362.96),
    @(253, 6, 2478, 1279.94),
    @(254, 4, 2135, 1063.66),
    @(255, 5, 2284, 1169.03),
    @(256, 6, 2558, 1327.72),
    @(257, 4, 2218, 1143.44)
)

foreach ($r in $results) {
    $row = $r[0]
    $ws.Range("H$row").Value = $r[1]
    $ws.Range("I$row").Value = $r[2]
    $ws.Range("J$row").Value = $r[3]
    $ws.Range("K$row").Value = "T"
    $ws.Range("L$row").Value = "T"
}

# Leave the sheet scrolled to the newly-entered block, with the last
# filled-in row range selected (matches the saved view after data entry).
[void]$excel.Goto($ws.Range("C237"))
$ws.Range("H254:L257").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Carry forward the saved selections on the other results sheets
#    that were also open/reviewed in this session.
# ---------------------------------------------------------------------
$wsReduced = $wb.Worksheets.Item("reduced_model_set_results")
$wsReduced.Range("L2").Select() | Out-Null

$wsNoCapacity = $wb.Worksheets.Item("no_capacity_results")
$wsNoCapacity.Range("K2").Select() | Out-Null

# Restore the originally active sheet/tab.
$ws.Activate() | Out-Null
